$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13: Week number, Total time, Comprehension scores text
$ws.Range("A13").Value = 12
$ws.Range("B13").Value = 0.87107638888888894
$ws.Range("B13").NumberFormat = "h:mm:ss"

$ws.Range("C13").Value = "Pirates of the Caribbean 2 (Audiovisual, English, Familiar):38; The Lion King (Audiovisual, English, Familiar):35; Harry Potter book 4 (Text-only, English, Familiar):33; [¿Nos Estamos Volviendo Más Estúpidos?](https://youtu.be/jQNvBfNjCeM) (Audiovisual, Spanish, New):37; Madre solo hay dos (Audiovisual, Spanish, New):32; W.I.T.C.H.  (Audiovisual, English, Familiar):37;"

# Update selection to match post-edit state (Excel moves active cell down after entry)
$ws.Range("C14").Select()
